$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = 13.25960146226672
$ws.Range("D10").Value = 12.12331036919332
$ws.Range("D11").Value = 13.84082545764194
$ws.Range("D12").Value = 7.210883324273421
$ws.Range("D20").Value = 14.10113279450843
$ws.Range("D21").Value = 12.14707767763526
$ws.Range("D22").Value = 13.50917967707238
